$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6049340963363647
$ws.Range("B1").Value = 0.7593464851379395
$ws.Range("C1").Value = 1.104142069816589
$ws.Range("D1").Value = 4.315260887145996
$ws.Range("E1").Value = 6.281651973724365
